$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 314.875
$ws.Range("I32").Value = 470.5
$ws.Range("J32").Value = 263
$ws.Range("K32").Value = 470.5
$ws.Range("L32").Value = 263
$ws.Range("M32").Value = -144.5
$ws.Range("N32").Value = -915

# Row 76
$ws.Range("H76").Value = 3266.6667
$ws.Range("I76").Value = 3266.6667
$ws.Range("K76").Value = 3266.6667
$ws.Range("M76").Value = -2951.6667

# Row 79
$ws.Range("H79").Value = 3266.6667
$ws.Range("I79").Value = 3266.6667
$ws.Range("K79").Value = 3266.6667
$ws.Range("M79").Value = -2174.6667

# Row 92
$ws.Range("H92").Value = 818.1111
$ws.Range("I92").Value = 577.8182
$ws.Range("J92").Value = 1195.7142
$ws.Range("K92").Value = 577.8182
$ws.Range("L92").Value = 1195.7142
$ws.Range("M92").Value = 670.1818
$ws.Range("N92").Value = -3691.7142

# Row 141
$ws.Range("H141").Value = 835499.2
$ws.Range("I141").Value = 910544.5600000001
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 2731633.68
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -2726453.68
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 21993.297
$ws.Range("I32").Value = 8668.929
$ws.Range("J32").Value = 35317.668
$ws.Range("K32").Value = 8668.929
$ws.Range("L32").Value = 35317.668
$ws.Range("M32").Value = -8381.929
$ws.Range("N32").Value = -35891.668

# Row 45
$ws.Range("H45").Value = 12988038
$ws.Range("I45").Value = 12988038
$ws.Range("K45").Value = 12988038
$ws.Range("M45").Value = -12987661

# Row 74
$ws.Range("H74").Value = 63375.832
$ws.Range("I74").Value = 501556
$ws.Range("J74").Value = 8603.3125
$ws.Range("K74").Value = 501556
$ws.Range("L74").Value = 8603.3125
$ws.Range("M74").Value = -500682
$ws.Range("N74").Value = -10351.3125

# Row 77
$ws.Range("H77").Value = 63375.832
$ws.Range("I77").Value = 501556
$ws.Range("J77").Value = 8603.3125
$ws.Range("K77").Value = 2507780
$ws.Range("L77").Value = 43016.5625
$ws.Range("M77").Value = -2503412
$ws.Range("N77").Value = -51752.5625

# Row 122
$ws.Range("H122").Value = 1534.325
$ws.Range("I122").Value = 1243.3928
$ws.Range("J122").Value = 2213.1667
$ws.Range("K122").Value = 3730.1784
$ws.Range("L122").Value = 6639.500100000001
$ws.Range("M122").Value = -1280.1784
$ws.Range("N122").Value = -11539.5001

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2431.45
$ws.Range("I105").Value = 2378.8572
$ws.Range("J105").Value = 2799.6
$ws.Range("K105").Value = 2378.8572
$ws.Range("L105").Value = 2799.6
$ws.Range("M105").Value = -631.8571999999999
$ws.Range("N105").Value = -6293.6

# Row 107
$ws.Range("H107").Value = 1435
$ws.Range("I107").Value = 1349.5333
$ws.Range("J107").Value = 1595.25
$ws.Range("K107").Value = 1349.5333
$ws.Range("L107").Value = 1595.25
$ws.Range("M107").Value = 570.4666999999999
$ws.Range("N107").Value = -5435.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2143.6216
$ws.Range("I31").Value = 1114.4878
$ws.Range("J31").Value = 3422.2424
$ws.Range("K31").Value = 1114.4878
$ws.Range("L31").Value = 3422.2424
$ws.Range("M31").Value = -819.4878000000001
$ws.Range("N31").Value = -4012.2424

# Row 34
$ws.Range("H34").Value = 2143.6216
$ws.Range("I34").Value = 1114.4878
$ws.Range("J34").Value = 3422.2424
$ws.Range("K34").Value = 1114.4878
$ws.Range("L34").Value = 3422.2424
$ws.Range("M34").Value = -912.4878000000001
$ws.Range("N34").Value = -3826.2424

# Row 68
$ws.Range("H68").Value = 20000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 20000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 20000
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -21498

# Row 71
$ws.Range("H71").Value = 20000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 20000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 60000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -67488

# Row 99
$ws.Range("H99").Value = 3093.0908
$ws.Range("J99").Value = 1833.3334
$ws.Range("L99").Value = 1833.3334
$ws.Range("N99").Value = -4829.3334

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 122
$ws.Range("H122").Value = 942.1905
$ws.Range("I122").Value = 843
$ws.Range("J122").Value = 1140.5714
$ws.Range("K122").Value = 2529
$ws.Range("L122").Value = 3421.7142
$ws.Range("M122").Value = -79
$ws.Range("N122").Value = -8321.7142

# Row 126
$ws.Range("H126").Value = 3093.0908
$ws.Range("J126").Value = 1833.3334
$ws.Range("L126").Value = 5500.0002
$ws.Range("N126").Value = -10440.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 915.3
$ws.Range("I98").Value = 933.125
$ws.Range("J98").Value = 844
$ws.Range("K98").Value = 2799.375
$ws.Range("L98").Value = 2532
$ws.Range("M98").Value = -1301.375
$ws.Range("N98").Value = -5528

# Row 113
$ws.Range("H113").Value = 1595440.5
$ws.Range("I113").Value = 2331541.5
$ws.Range("J113").Value = 554.6667
$ws.Range("K113").Value = 6994624.5
$ws.Range("L113").Value = 1664.0001
$ws.Range("M113").Value = -6992454.5
$ws.Range("N113").Value = -6004.0001

# Row 122
$ws.Range("H122").Value = 37377.8
$ws.Range("I122").Value = 42858.848
$ws.Range("J122").Value = 1751
$ws.Range("K122").Value = 385729.632
$ws.Range("L122").Value = 15759
$ws.Range("M122").Value = -383279.632
$ws.Range("N122").Value = -20659

# Row 131
$ws.Range("H131").Value = 2180.9375
$ws.Range("I131").Value = 3783.3333
$ws.Range("J131").Value = 1811.1538
$ws.Range("K131").Value = 11349.9999
$ws.Range("L131").Value = 5433.4614
$ws.Range("M131").Value = -6309.999899999999
$ws.Range("N131").Value = -15513.4614

# Row 137
$ws.Range("H137").Value = 9814437
$ws.Range("I137").Value = 2890.0625
$ws.Range("J137").Value = 18076792
$ws.Range("K137").Value = 8670.1875
$ws.Range("L137").Value = 54230376
$ws.Range("M137").Value = -3570.1875
$ws.Range("N137").Value = -54240576

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5029.375
$ws.Range("I70").Value = 4624.625
$ws.Range("J70").Value = 5838.875
$ws.Range("K70").Value = 4624.625
$ws.Range("L70").Value = 5838.875
$ws.Range("M70").Value = -4354.625
$ws.Range("N70").Value = -6378.875

# Row 73
$ws.Range("H73").Value = 5029.375
$ws.Range("I73").Value = 4624.625
$ws.Range("J73").Value = 5838.875
$ws.Range("K73").Value = 4624.625
$ws.Range("L73").Value = 5838.875
$ws.Range("M73").Value = -3688.625
$ws.Range("N73").Value = -7710.875

# Row 80
$ws.Range("H80").Value = 2523.3333
$ws.Range("I80").Value = 2410
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 2410
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -1412
$ws.Range("N80").Value = -4746

# Row 83
$ws.Range("H83").Value = 2523.3333
$ws.Range("I83").Value = 2410
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 12050
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -7058
$ws.Range("N83").Value = -23734

# Row 126
$ws.Range("H126").Value = 20835058
$ws.Range("I126").Value = 66667784
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 200003352
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -200000882
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2856.6667
$ws.Range("I61").Value = 2340
$ws.Range("J61").Value = 3502.5
$ws.Range("K61").Value = 2340
$ws.Range("L61").Value = 3502.5
$ws.Range("M61").Value = -2138
$ws.Range("N61").Value = -3906.5

# Row 113
$ws.Range("H113").Value = 2856.6667
$ws.Range("I113").Value = 2340
$ws.Range("J113").Value = 3502.5
$ws.Range("K113").Value = 2340
$ws.Range("L113").Value = 3502.5
$ws.Range("M113").Value = -170
$ws.Range("N113").Value = -7842.5

# Row 122
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3675
$ws.Range("J122").Value = 3383.3333
$ws.Range("K122").Value = 11025
$ws.Range("L122").Value = 10149.9999
$ws.Range("M122").Value = -8575
$ws.Range("N122").Value = -15049.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 42323.414
$ws.Range("I122").Value = 661.64703
$ws.Range("J122").Value = 101344.25
$ws.Range("K122").Value = 1984.94109
$ws.Range("L122").Value = 304032.75
$ws.Range("M122").Value = 465.0589100000002
$ws.Range("N122").Value = -308932.75
